$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C, E, F, G columns for rows with new S-curve-based projections (PLF ramp and future forecasts)

$ws.Range("C25").Value = 1.605796365659359
$efg25 = New-Object 'object[,]' 1,3
$efg25[0,0] = 2.418585580357664
$efg25[0,1] = 0.1771200564062623
$efg25[0,2] = 0.117597138249909
$ws.Range("E25:G25").Value = $efg25

$ws.Range("C26").Value = 1.595170373197115
$efg26 = New-Object 'object[,]' 1,3
$efg26[0,0] = 2.423189895590971
$efg26[0,1] = 0.1774572438022588
$efg26[0,2] = 0.1168189659166369
$ws.Range("E26:G26").Value = $efg26

$ws.Range("C31").Value = 1.559211995598527
$efg31 = New-Object 'object[,]' 1,3
$efg31[0,0] = 2.259060943174633
$efg31[0,1] = 0.1654376032545038
$efg31[0,2] = 0.1141856293416303
$ws.Range("E31:G31").Value = $efg31

$ws.Range("C32").Value = 1.567714000147012
$efg32 = New-Object 'object[,]' 1,3
$efg32[0,0] = 2.289265470940818
$efg32[0,1] = 0.1676495686714481
$efg32[0,2] = 0.1148082558624463
$ws.Range("E32:G32").Value = $efg32

$ws.Range("C33").Value = 1.564079939769674
$efg33 = New-Object 'object[,]' 1,3
$efg33[0,0] = 2.264544154068898
$efg33[0,1] = 0.1658391547359842
$efg33[0,2] = 0.1145421230514987
$ws.Range("E33:G33").Value = $efg33

$ws.Range("C34").Value = 1.554941631335987
$efg34 = New-Object 'object[,]' 1,3
$efg34[0,0] = 2.187896682526145
$efg34[0,1] = 0.1602260374688904
$efg34[0,2] = 0.1138728981465057
$ws.Range("E34:G34").Value = $efg34

$ws.Range("C35").Value = 1.550489742549896
$efg35 = New-Object 'object[,]' 1,3
$efg35[0,0] = 2.239272596246528
$efg35[0,1] = 0.1639884450553657
$efg35[0,2] = 0.1135468733825648
$ws.Range("E35:G35").Value = $efg35

$ws.Range("C38").Value = 1.427985993580319
$efg38 = New-Object 'object[,]' 1,3
$efg38[0,0] = 1.945719472026729
$efg38[0,1] = 0.1424906959815635
$efg38[0,2] = 0.1045755675484081
$ws.Range("E38:G38").Value = $efg38

$ws.Range("C39").Value = 1.410548045123722
$efg39 = New-Object 'object[,]' 1,3
$efg39[0,0] = 1.882167269668976
$efg39[0,1] = 0.1378365833639389
$efg39[0,2] = 0.1032985358653757
$ws.Range("E39:G39").Value = $efg39

$ws.Range("C40").Value = 1.405378496021678
$efg40 = New-Object 'object[,]' 1,3
$efg40[0,0] = 1.853692822352101
$efg40[0,1] = 0.1357513167701655
$efg40[0,2] = 0.1029199547492121
$ws.Range("E40:G40").Value = $efg40

$ws.Range("C41").Value = 1.408893757585069
$efg41 = New-Object 'object[,]' 1,3
$efg41[0,0] = 1.834384170021593
$efg41[0,1] = 0.1343372880015813
$efg41[0,2] = 0.1031773875774929
$ws.Range("E41:G41").Value = $efg41

$ws.Range("C42").Value = 1.380351075228878
$efg42 = New-Object 'object[,]' 1,3
$efg42[0,0] = 1.815920663940857
$efg42[0,1] = 0.1329851517509411
$efg42[0,2] = 0.1010871239333315
$ws.Range("E42:G42").Value = $efg42

$ws.Range("C43").Value = 1.363977034209024
$efg43 = New-Object 'object[,]' 1,3
$efg43[0,0] = 1.778751645324429
$efg43[0,1] = 0.1302631564130984
$efg43[0,2] = 0.09988800528622285
$ws.Range("E43:G43").Value = $efg43

$ws.Range("C44").Value = 1.380901000651561
$efg44 = New-Object 'object[,]' 1,3
$efg44[0,0] = 1.766531543861792
$efg44[0,1] = 0.1293682428413271
$efg44[0,2] = 0.1011273965714701
$ws.Range("E44:G44").Value = $efg44

$ws.Range("C47").Value = 1.342218401785662
$efg47 = New-Object 'object[,]' 1,3
$efg47[0,0] = 1.688706241602388
$efg47[0,1] = 0.1236688696051807
$efg47[0,2] = 0.0982945573497727
$ws.Range("E47:G47").Value = $efg47

$ws.Range("C48").Value = 1.331246744782096
$efg48 = New-Object 'object[,]' 1,3
$efg48[0,0] = 1.669934466960038
$efg48[0,1] = 0.122294157951187
$efg48[0,2] = 0.0974910709967885
$ws.Range("E48:G48").Value = $efg48

$ws.Range("C49").Value = 1.315140810373198
$efg49 = New-Object 'object[,]' 1,3
$efg49[0,0] = 1.639036470100827
$efg49[0,1] = 0.1200314077756342
$efg49[0,2] = 0.09631158657657721
$ws.Range("E49:G49").Value = $efg49

$ws.Range("C50").Value = 1.304783119122762
$efg50 = New-Object 'object[,]' 1,3
$efg50[0,0] = 1.625178676835295
$efg50[0,1] = 0.1190165612699776
$efg50[0,2] = 0.09555306272138883
$ws.Range("E50:G50").Value = $efg50

$ws.Range("C51").Value = 1.302501009379456
$efg51 = New-Object 'object[,]' 1,3
$efg51[0,0] = 1.598601300647098
$efg51[0,1] = 0.1170702227125106
$efg51[0,2] = 0.09538593718746423
$ws.Range("E51:G51").Value = $efg51

$ws.Range("C52").Value = 1.296052864122837
$efg52 = New-Object 'object[,]' 1,3
$efg52[0,0] = 1.588232197226381
$efg52[0,1] = 0.1163108631109007
$efg52[0,2] = 0.09491372075615678
$ws.Range("E52:G52").Value = $efg52

$ws.Range("C53").Value = 1.276281382275522
$efg53 = New-Object 'object[,]' 1,3
$efg53[0,0] = 1.549199684641708
$efg53[0,1] = 0.1134523986898679
$efg53[0,2] = 0.09346579763593625
$ws.Range("E53:G53").Value = $efg53

$ws.Range("C54").Value = 1.190921954421548
$efg54 = New-Object 'object[,]' 1,3
$efg54[0,0] = 1.825434111347368
$efg54[0,1] = 0.1336818491739903
$efg54[0,2] = 0.08721467846980514
$ws.Range("E54:G54").Value = $efg54

$ws.Range("C55").Value = 1.150148931499757
$efg55 = New-Object 'object[,]' 1,3
$efg55[0,0] = 1.699512378614168
$efg55[0,1] = 0.1244602344477583
$efg55[0,2] = 0.08422875141458244
$ws.Range("E55:G55").Value = $efg55

$ws.Range("C56").Value = 1.142005350412699
$efg56 = New-Object 'object[,]' 1,3
$efg56[0,0] = 1.681497877118
$efg56[0,1] = 0.1231409801087573
$efg56[0,2] = 0.08363237328630661
$ws.Range("E56:G56").Value = $efg56

$ws.Range("C57").Value = 1.135647554608852
$efg57 = New-Object 'object[,]' 1,3
$efg57[0,0] = 1.667433730055871
$efg57[0,1] = 0.1221110217144045
$efg57[0,2] = 0.08316677340820329
$ws.Range("E57:G57").Value = $efg57

$ws.Range("C58").Value = 1.124752312390587
$efg58 = New-Object 'object[,]' 1,3
$efg58[0,0] = 1.643332249989076
$efg58[0,1] = 0.1203460002309498
$efg58[0,2] = 0.08236888313219592
$ws.Range("E58:G58").Value = $efg58

$ws.Range("C59").Value = 1.106894511493559
$efg59 = New-Object 'object[,]' 1,3
$efg59[0,0] = 1.603828820534797
$efg59[0,1] = 0.1174530491979133
$efg59[0,2] = 0.08106110443382726
$ws.Range("E59:G59").Value = $efg59

$ws.Range("C60").Value = 1.079657624613052
$efg60 = New-Object 'object[,]' 1,3
$efg60[0,0] = 1.543577816189835
$efg60[0,1] = 0.1130406929121641
$efg60[0,2] = 0.07906646799020264
$ws.Range("E60:G60").Value = $efg60

$ws.Range("C61").Value = 1.042361339290171
$efg61 = New-Object 'object[,]' 1,3
$efg61[0,0] = 1.461074310598378
$efg61[0,1] = 0.1069987212396495
$efg61[0,2] = 0.07633515254129643
$ws.Range("E61:G61").Value = $efg61

$ws.Range("C62").Value = 0.9981777361154145
$efg62 = New-Object 'object[,]' 1,3
$efg62[0,0] = 1.363335296645078
$efg62[0,1] = 0.09984100897794779
$efg62[0,2] = 0.07309945877461668
$ws.Range("E62:G62").Value = $efg62

$ws.Range("C63").Value = 0.9539941329406577
$efg63 = New-Object 'object[,]' 1,3
$efg63[0,0] = 1.265596282691778
$efg63[0,1] = 0.09268329671624609
$efg63[0,2] = 0.06986376500793694
$ws.Range("E63:G63").Value = $efg63

$ws.Range("C64").Value = 0.9166978476177772
$efg64 = New-Object 'object[,]' 1,3
$efg64[0,0] = 1.18309277710032
$efg64[0,1] = 0.08664132504373147
$efg64[0,2] = 0.06713244955903074
$ws.Range("E64:G64").Value = $efg64

$ws.Range("C65").Value = 0.8894609607372701
$efg65 = New-Object 'object[,]' 1,3
$efg65[0,0] = 1.122841772755358
$efg65[0,1] = 0.08222896875798222
$efg65[0,2] = 0.06513781311540612
$ws.Range("E65:G65").Value = $efg65

$ws.Range("C66").Value = 0.8716031598402422
$efg66 = New-Object 'object[,]' 1,3
$efg66[0,0] = 1.083338343301079
$efg66[0,1] = 0.07933601772494581
$efg66[0,2] = 0.06383003441703744
$ws.Range("E66:G66").Value = $efg66

$ws.Range("C67").Value = 0.8607079176219773
$efg67 = New-Object 'object[,]' 1,3
$efg67[0,0] = 1.059236863234285
$efg67[0,1] = 0.07757099624149107
$efg67[0,2] = 0.06303214414103009
$ws.Range("E67:G67").Value = $efg67

$ws.Range("C68").Value = 0.8543501218181296
$efg68 = New-Object 'object[,]' 1,3
$efg68[0,0] = 1.045172716172156
$efg68[0,1] = 0.07654103784713821
$efg68[0,2] = 0.06256654426292677
$ws.Range("E68:G68").Value = $efg68

$ws.Range("C69").Value = 0.8411293014866855
$efg69 = New-Object 'object[,]' 1,3
$efg69[0,0] = 1.020995265387932
$efg69[0,1] = 0.07477045280708841
$efg69[0,2] = 0.06159834513784302
$ws.Range("E69:G69").Value = $efg69

$ws.Range("C70").Value = 0.8354032930588543
$efg70 = New-Object 'object[,]' 1,3
$efg70[0,0] = 1.01404481498268
$efg70[0,1] = 0.07426145110881276
$efg70[0,2] = 0.06117901288681298
$ws.Range("E70:G70").Value = $efg70

$ws.Range("C71").Value = 0.8238011807776481
$efg71 = New-Object 'object[,]' 1,3
$efg71[0,0] = 0.9999617225417519
$efg71[0,1] = 0.0732301052892685
$efg71[0,2] = 0.06032935646019393
$ws.Range("E71:G71").Value = $efg71

$ws.Range("C72").Value = 0.8020049483276099
$efg72 = New-Object 'object[,]' 1,3
$efg72[0,0] = 0.9735046129208524
$efg72[0,1] = 0.07129257420231504
$efg72[0,2] = 0.05873315496443215
$ws.Range("E72:G72").Value = $efg72

$ws.Range("C73").Value = 0.7663231432013164
$efg73 = New-Object 'object[,]' 1,3
$efg73[0,0] = 0.9301926583497191
$efg73[0,1] = 0.06812071379803274
$efg73[0,2] = 0.05612007259597127
$ws.Range("E73:G73").Value = $efg73

$ws.Range("C74").Value = 0.7192755596214111
$efg74 = New-Object 'object[,]' 1,3
$efg74[0,0] = 0.8730844824745896
$efg74[0,1] = 0.06393851598191665
$efg74[0,2] = 0.05267464121445328
$ws.Range("E74:G74").Value = $efg74

$ws.Range("C75").Value = 0.6722279760415057
$efg75 = New-Object 'object[,]' 1,3
$efg75[0,0] = 0.8159763065994602
$efg75[0,1] = 0.05975631816580056
$efg75[0,2] = 0.0492292098329353
$ws.Range("E75:G75").Value = $efg75

$ws.Range("C76").Value = 0.6365461709152123
$efg76 = New-Object 'object[,]' 1,3
$efg76[0,0] = 0.772664352028327
$efg76[0,1] = 0.05658445776151827
$efg76[0,2] = 0.04661612746447441
$ws.Range("E76:G76").Value = $efg76

$ws.Range("C77").Value = 0.6147499384651741
$efg77 = New-Object 'object[,]' 1,3
$efg77[0,0] = 0.7462072424074275
$efg77[0,1] = 0.0546469266745648
$efg77[0,2] = 0.04501992596871264
$ws.Range("E77:G77").Value = $efg77

$ws.Range("C78").Value = 0.6031478261839679
$efg78 = New-Object 'object[,]' 1,3
$efg78[0,0] = 0.7321241499664989
$efg78[0,1] = 0.05361558085502055
$efg78[0,2] = 0.0441702695420936
$ws.Range("E78:G78").Value = $efg78

$ws.Range("C79").Value = 0.5974218177561367
$efg79 = New-Object 'object[,]' 1,3
$efg79[0,0] = 0.7251736995612473
$efg79[0,1] = 0.05310657915674488
$efg79[0,2] = 0.04375093729106354
$ws.Range("E79:G79").Value = $efg79

$ws.Range("C80").Value = 0.5923445785117503
$efg80 = New-Object 'object[,]' 1,3
$efg80[0,0] = 0.7190107502731914
$efg80[0,1] = 0.05265524845569606
$efg80[0,2] = 0.04337911629425564
$ws.Range("E80:G80").Value = $efg80

$ws.Range("C81").Value = 0.5923445785117503
$efg81 = New-Object 'object[,]' 1,3
$efg81[0,0] = 0.7190107502731914
$efg81[0,1] = 0.05265524845569606
$efg81[0,2] = 0.04337911629425564
$ws.Range("E81:G81").Value = $efg81

$ws.Range("C82").Value = 0.5923445785117503
$efg82 = New-Object 'object[,]' 1,3
$efg82[0,0] = 0.7190107502731914
$efg82[0,1] = 0.05265524845569606
$efg82[0,2] = 0.04337911629425564
$ws.Range("E82:G82").Value = $efg82

$ws.Range("C83").Value = 0.5923445785117503
$efg83 = New-Object 'object[,]' 1,3
$efg83[0,0] = 0.7190107502731914
$efg83[0,1] = 0.05265524845569606
$efg83[0,2] = 0.04337911629425564
$ws.Range("E83:G83").Value = $efg83

$ws.Range("C84").Value = 0.5923445785117503
$efg84 = New-Object 'object[,]' 1,3
$efg84[0,0] = 0.7190107502731914
$efg84[0,1] = 0.05265524845569606
$efg84[0,2] = 0.04337911629425564
$ws.Range("E84:G84").Value = $efg84
